# Commit: "mudando de DRH pra DP" (changing from DRH to DP)
#
# The institution name "DIRETORIA DE RECURSOS HUMANOS" (DRH) was renamed to
# "DIRETORIA DE PESSOAL" (DP), and every other reference to the "DRH"
# acronym in the body of the note was updated to "DP" accordingly.

$d = $word.ActiveDocument

# 1) Full directorate name in the letterhead.
$d.Content.Find.Execute("DIRETORIA DE RECURSOS HUMANOS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "DIRETORIA DE PESSOAL", 2) | Out-Null

# 2) Acronym used in the note number ("{nota_bg}/DRH-1/2025.") and in the
#    "À DRH para providências;" line. Both simply swap the "DRH" acronym
#    for "DP", so a single targeted replace covers both spots.
$d.Content.Find.Execute("DRH", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "DP", 2) | Out-Null

# 3) Drop the stray "_GoBack" bookmark left over from the previous save.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
